$d = $word.ActiveDocument

function Get-ContainingParagraph($doc, $pos) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs($i)
        if ($pos -ge $p.Range.Start -and $pos -lt $p.Range.End) {
            return $p
        }
    }
    return $null
}

# ------------------------------------------------------------------
# Change 1: the run that reads " Marek Switajski & Michael Werner"
# (directly after "Gestaltung & Produktion:") becomes " Michael
# Werner". It must stay its own <w:r>, separate from the preceding
# "Gestaltung & Produktion:" run rather than merging into it.
# ------------------------------------------------------------------

$oldText1 = " Marek Switajski & Michael Werner"
$newText1 = " Michael Werner"

$search1 = $d.Content.Duplicate
$found1 = $search1.Find.Execute($oldText1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find first target run text"
}
$run1Start = $search1.Start
$run1End = $search1.End

# Drop the old run content entirely (the now-empty run disappears).
$oldRange1 = $d.Range($run1Start, $run1End)
$oldRange1.Text = ""

# Inserting new text right at the end of a paragraph (just before the
# paragraph mark) always creates a brand-new run there instead of
# being folded into the preceding run - so stage the replacement text
# there first.
$para1 = Get-ContainingParagraph $d $run1Start
$paraEnd1 = $para1.Range.End - 1
$tempIns1 = $d.Range($paraEnd1, $paraEnd1)
$tempIns1.InsertAfter($newText1)

# Relocate that freshly-created, still-standalone run from the end of
# the paragraph to where the old text used to live, via Cut/Paste, so
# it keeps its own run identity instead of merging with its new
# neighbours.
$newRunRange1 = $d.Range($paraEnd1, $paraEnd1 + $newText1.Length)
$newRunRange1.Cut()
$destination1 = $d.Range($run1Start, $run1Start)
$destination1.Paste()

# ------------------------------------------------------------------
# Change 2: the run "Marek Switajski & Michael Werner" that follows
# "Web-Administration und -Programmierung: " is split into three
# runs: "Michael Werner", " & ", "Marek Switajski".
# ------------------------------------------------------------------

$oldText2 = "Marek Switajski & Michael Werner"

$search2 = $d.Content.Duplicate
$found2 = $search2.Find.Execute($oldText2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find second target run text"
}
$run2Start = $search2.Start
$run2End = $search2.End

$oldRange2 = $d.Range($run2Start, $run2End)
$oldRange2.Text = ""

$cursor = $run2Start
$pieces = @("Michael Werner", " & ", "Marek Switajski")
foreach ($piece in $pieces) {
    $insertionPoint = $d.Range($cursor, $cursor)
    $insertionPoint.InsertAfter($piece)
    $cursor = $insertionPoint.End
}

Write-Output "Done."
